$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$textCells = @('D5', 'D6', 'D10', 'D12', 'D14', 'D18', 'D20', 'D21', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D35', 'D36', 'D38', 'D39', 'D40', 'D41', 'D43', 'D44', 'D46', 'D47', 'D48', 'D50')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values from the refreshed crypto feed
$ws.Range('D2').Value = '66.977.77'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '3.102.10'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '576.97'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '178.51'
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.099.84'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').Value = '6.34'
$ws.Range('E10').Value = '  -2.91%  '
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('D12').Value = '0.468'
$ws.Range('E12').Value = '  -2.49%  '
$ws.Range('E13').Value = '  -3.01%  '
$ws.Range('D14').Value = '36.17'
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '3.620.30'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '66.943.86'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').Value = '7.01'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.100.99'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '16.68'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '481.29'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('E23').Value = '  -2.09%  '
$ws.Range('D24').Value = '83.70'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '12.65'
$ws.Range('E25').Value = '  -3.98%  '
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  -2.23%  '
$ws.Range('D27').Value = '10.12'
$ws.Range('E27').Value = '  -4.45%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '8.00'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('D30').Value = '2.28'
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').Value = '2.61'
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').Value = '28.02'
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('E33').Value = '  -1.97%  '
$ws.Range('D34').Value = '0.0₃0942'
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').Value = '48.50'
$ws.Range('E36').Value = '  +3.04%  '
$ws.Range('E37').Value = '  -4.72%  '
$ws.Range('D38').Value = '0.937'
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('D39').Value = '0.312'
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('D40').Value = '49.08'
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('D41').Value = '2.00'
$ws.Range('E41').Value = '  -2.25%  '
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').Value = '8.34'
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('D44').Value = '2.68'
$ws.Range('E44').Value = '  +3.67%  '
$ws.Range('D45').Value = '2.798.97'
$ws.Range('E45').Value = '  -0.15%  '
$ws.Range('D46').Value = '373.83'
$ws.Range('E46').Value = '  -4.64%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0344'
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '135.60'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '25.12'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('E51').Value = '  +1.60%  '
